$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 31/12/2025 18:55:54"
$ws1.Range("A3").Value = "Total filas: 1256"

$ws1.Cells.Item(1240,2).Value = "18:55:43"
$ws1.Cells.Item(1240,3).Value = "19:01"
$ws1.Cells.Item(1240,4).Value = "17_ROMERO"
$ws1.Cells.Item(1240,5).Value = 6
$ws1.Cells.Item(1240,6).Value = "LP1912"
$ws1.Cells.Item(1240,7).Value = "31/12/2025"

$ws1.Cells.Item(1241,2).Value = "18:55:43"
$ws1.Cells.Item(1241,3).Value = "19:03"
$ws1.Cells.Item(1241,4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(1241,5).Value = 8
$ws1.Cells.Item(1241,6).Value = "LP1912"
$ws1.Cells.Item(1241,7).Value = "31/12/2025"

$ws1.Cells.Item(1242,2).Value = "18:55:43"
$ws1.Cells.Item(1242,3).Value = "19:05"
$ws1.Cells.Item(1242,4).Value = "16_SANTA ANA"
$ws1.Cells.Item(1242,5).Value = 10
$ws1.Cells.Item(1242,6).Value = "LP1912"
$ws1.Cells.Item(1242,7).Value = "31/12/2025"

$ws1.Cells.Item(1243,2).Value = "18:55:43"
$ws1.Cells.Item(1243,3).Value = "19:11"
$ws1.Cells.Item(1243,4).Value = "81_EL PELIGRO"
$ws1.Cells.Item(1243,5).Value = 16
$ws1.Cells.Item(1243,6).Value = "LP1912"
$ws1.Cells.Item(1243,7).Value = "31/12/2025"

$ws1.Cells.Item(1244,2).Value = "18:55:43"
$ws1.Cells.Item(1244,3).Value = "19:14"
$ws1.Cells.Item(1244,4).Value = "14_ABASTO"
$ws1.Cells.Item(1244,5).Value = 19
$ws1.Cells.Item(1244,6).Value = "LP1912"
$ws1.Cells.Item(1244,7).Value = "31/12/2025"

$ws1.Cells.Item(1245,2).Value = "18:55:43"
$ws1.Cells.Item(1245,3).Value = "19:17"
$ws1.Cells.Item(1245,4).Value = "16_SANTA ANA"
$ws1.Cells.Item(1245,5).Value = 22
$ws1.Cells.Item(1245,6).Value = "LP1912"
$ws1.Cells.Item(1245,7).Value = "31/12/2025"

$ws1.Cells.Item(1246,2).Value = "18:55:43"
$ws1.Cells.Item(1246,3).Value = "19:20"
$ws1.Cells.Item(1246,4).Value = "215C_EL PATO"
$ws1.Cells.Item(1246,5).Value = 25
$ws1.Cells.Item(1246,6).Value = "LP1912"
$ws1.Cells.Item(1246,7).Value = "31/12/2025"

$ws1.Cells.Item(1247,2).Value = "18:55:43"
$ws1.Cells.Item(1247,3).Value = "19:29"
$ws1.Cells.Item(1247,4).Value = "16_SANTA ANA"
$ws1.Cells.Item(1247,5).Value = 34
$ws1.Cells.Item(1247,6).Value = "LP1912"
$ws1.Cells.Item(1247,7).Value = "31/12/2025"

$ws1.Cells.Item(1248,2).Value = "18:55:43"
$ws1.Cells.Item(1248,3).Value = "19:29"
$ws1.Cells.Item(1248,4).Value = "225_GOMEZ"
$ws1.Cells.Item(1248,5).Value = 34
$ws1.Cells.Item(1248,6).Value = "LP1912"
$ws1.Cells.Item(1248,7).Value = "31/12/2025"

$ws1.Cells.Item(1249,2).Value = "18:55:43"
$ws1.Cells.Item(1249,3).Value = "19:33"
$ws1.Cells.Item(1249,4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(1249,5).Value = 38
$ws1.Cells.Item(1249,6).Value = "LP1912"
$ws1.Cells.Item(1249,7).Value = "31/12/2025"

$ws1.Cells.Item(1250,2).Value = "18:55:43"
$ws1.Cells.Item(1250,3).Value = "19:40"
$ws1.Cells.Item(1250,4).Value = "17X38_ROMERO"
$ws1.Cells.Item(1250,5).Value = 45
$ws1.Cells.Item(1250,6).Value = "LP1912"
$ws1.Cells.Item(1250,7).Value = "31/12/2025"

$ws1.Cells.Item(1251,2).Value = "18:55:43"
$ws1.Cells.Item(1251,3).Value = "19:44"
$ws1.Cells.Item(1251,4).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(1251,5).Value = 49
$ws1.Cells.Item(1251,6).Value = "LP1912"
$ws1.Cells.Item(1251,7).Value = "31/12/2025"

$ws1.Cells.Item(1252,2).Value = "18:55:43"
$ws1.Cells.Item(1252,3).Value = "19:51"
$ws1.Cells.Item(1252,4).Value = "81_EL PELIGRO"
$ws1.Cells.Item(1252,5).Value = 56
$ws1.Cells.Item(1252,6).Value = "LP1912"
$ws1.Cells.Item(1252,7).Value = "31/12/2025"

$ws1.Cells.Item(1253,2).Value = "18:55:43"
$ws1.Cells.Item(1253,3).Value = "19:58"
$ws1.Cells.Item(1253,4).Value = "14X44_ABASTO"
$ws1.Cells.Item(1253,5).Value = 63
$ws1.Cells.Item(1253,6).Value = "LP1912"
$ws1.Cells.Item(1253,7).Value = "31/12/2025"

$ws1.Cells.Item(1254,2).Value = "18:55:43"
$ws1.Cells.Item(1254,3).Value = "20:00"
$ws1.Cells.Item(1254,4).Value = "215C_EL PATO"
$ws1.Cells.Item(1254,5).Value = 65
$ws1.Cells.Item(1254,6).Value = "LP1912"
$ws1.Cells.Item(1254,7).Value = "31/12/2025"

$ws1.Cells.Item(1255,2).Value = "18:55:43"
$ws1.Cells.Item(1255,3).Value = "20:10"
$ws1.Cells.Item(1255,4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(1255,5).Value = 75
$ws1.Cells.Item(1255,6).Value = "LP1912"
$ws1.Cells.Item(1255,7).Value = "31/12/2025"

$ws1.Cells.Item(1256,2).Value = "18:55:43"
$ws1.Cells.Item(1256,3).Value = "20:26"
$ws1.Cells.Item(1256,4).Value = "15_ABASTO"
$ws1.Cells.Item(1256,5).Value = 91
$ws1.Cells.Item(1256,6).Value = "LP1912"
$ws1.Cells.Item(1256,7).Value = "31/12/2025"

$ws1.Cells.Item(1257,2).Value = "18:55:43"
$ws1.Cells.Item(1257,3).Value = "20:28"
$ws1.Cells.Item(1257,4).Value = "10_OLMOS"
$ws1.Cells.Item(1257,5).Value = 93
$ws1.Cells.Item(1257,6).Value = "LP1912"
$ws1.Cells.Item(1257,7).Value = "31/12/2025"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 31/12/2025 18:55:54"
$ws2.Range("A3").Value = "Total filas: 88"

$ws2.Cells.Item(88,2).Value = "31/12/2025"
$ws2.Cells.Item(88,3).Value = "18:55:43"
$ws2.Cells.Item(88,4).Value = "19:20"
$ws2.Cells.Item(88,5).Value = "215C_EL PATO"
$ws2.Cells.Item(88,6).Value = 25
$ws2.Cells.Item(88,7).Value = "LP1912"

$ws2.Cells.Item(89,2).Value = "31/12/2025"
$ws2.Cells.Item(89,3).Value = "18:55:43"
$ws2.Cells.Item(89,4).Value = "20:00"
$ws2.Cells.Item(89,5).Value = "215C_EL PATO"
$ws2.Cells.Item(89,6).Value = 65
$ws2.Cells.Item(89,7).Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 31/12/2025 18:55:54"
$ws3.Range("A3").Value = "Total filas: 146"

$ws3.Cells.Item(147,2).Value = "31/12/2025"
$ws3.Cells.Item(147,3).Value = "18:55:54"
$ws3.Cells.Item(147,4).Value = "19:10"
$ws3.Cells.Item(147,5).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(147,6).Value = 15
$ws3.Cells.Item(147,7).Value = "L6173"

